$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F ("想去人数") for the listed rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 263
$ws1.Range("F7").Value  = 61
$ws1.Range("F8").Value  = 257
$ws1.Range("F13").Value = 2163
$ws1.Range("F14").Value = 55
$ws1.Range("F16").Value = 504
$ws1.Range("F17").Value = 492
$ws1.Range("F18").Value = 153
$ws1.Range("F22").Value = 1652
$ws1.Range("F23").Value = 3822
$ws1.Range("F27").Value = 1142
$ws1.Range("F28").Value = 181
$ws1.Range("F29").Value = 2031
$ws1.Range("F30").Value = 567
$ws1.Range("F31").Value = 465
$ws1.Range("F32").Value = 81
$ws1.Range("F35").Value = 455
$ws1.Range("F36").Value = 659
$ws1.Range("F38").Value = 393

# Sheet "演出" (sheet2): update column F for its single data row
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 21

# Sheet "全部类型" (sheet4): update column F for the listed rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 263
$ws4.Range("F7").Value  = 61
$ws4.Range("F8").Value  = 257
$ws4.Range("F13").Value = 2163
$ws4.Range("F14").Value = 55
$ws4.Range("F15").Value = 21
$ws4.Range("F17").Value = 504
$ws4.Range("F18").Value = 492
$ws4.Range("F19").Value = 153
$ws4.Range("F23").Value = 1652
$ws4.Range("F24").Value = 3822
$ws4.Range("F28").Value = 1142
$ws4.Range("F29").Value = 182
$ws4.Range("F30").Value = 2031
$ws4.Range("F31").Value = 567
$ws4.Range("F32").Value = 465
$ws4.Range("F33").Value = 81
$ws4.Range("F36").Value = 455
$ws4.Range("F37").Value = 659
$ws4.Range("F39").Value = 393
